# Update computed market-board value columns (H-N) for several leve rows
# across multiple worksheets, reflecting refreshed price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 1007
$ws.Range("I33").Value = 325.81818
$ws.Range("K33").Value = 325.81818
$ws.Range("M33").Value = -96.81817999999998
# Row 88
$ws.Range("H88").Value = 4800.3
$ws.Range("I88").Value = 1099.5
$ws.Range("J88").Value = 5725.5
$ws.Range("K88").Value = 1099.5
$ws.Range("L88").Value = 5725.5
$ws.Range("M88").Value = -693.5
$ws.Range("N88").Value = -6537.5
# Row 91
$ws.Range("H91").Value = 4800.3
$ws.Range("I91").Value = 1099.5
$ws.Range("J91").Value = 5725.5
$ws.Range("K91").Value = 1099.5
$ws.Range("L91").Value = 5725.5
$ws.Range("M91").Value = 304.5
$ws.Range("N91").Value = -8533.5
# Row 96
$ws.Range("H96").Value = 700.5
$ws.Range("I96").Value = 497.5
$ws.Range("J96").Value = 781.7
$ws.Range("K96").Value = 1492.5
$ws.Range("L96").Value = 2345.1
$ws.Range("M96").Value = -119.5
$ws.Range("N96").Value = -5091.1
# Row 111
$ws.Range("H111").Value = 5564876
$ws.Range("I111").Value = 13901.7
$ws.Range("J111").Value = 12503594
$ws.Range("K111").Value = 41705.10000000001
$ws.Range("L111").Value = 37510782
$ws.Range("M111").Value = -38638.10000000001
$ws.Range("N111").Value = -37516916
# Row 112
$ws.Range("H112").Value = 959.2174
$ws.Range("J112").Value = 959.2174
$ws.Range("L112").Value = 2877.6522
$ws.Range("N112").Value = -5093.6522
# Row 129
$ws.Range("H129").Value = 959.975
$ws.Range("I129").Value = 617.6
$ws.Range("J129").Value = 1074.1
$ws.Range("K129").Value = 1852.8
$ws.Range("L129").Value = 3222.3
$ws.Range("M129").Value = 3147.2
$ws.Range("N129").Value = -13222.3
# Row 138
$ws.Range("H138").Value = 6190.7456
$ws.Range("I138").Value = 1212.8649
$ws.Range("K138").Value = 3638.5947
$ws.Range("M138").Value = 1501.4053

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 29482.184
$ws.Range("I32").Value = 4451.264
$ws.Range("K32").Value = 4451.264
$ws.Range("M32").Value = -4164.264
# Row 74
$ws.Range("H74").Value = 2083.8125
$ws.Range("I74").Value = 1742.7
$ws.Range("J74").Value = 2652.3333
$ws.Range("K74").Value = 1742.7
$ws.Range("L74").Value = 2652.3333
$ws.Range("M74").Value = -868.7
$ws.Range("N74").Value = -4400.3333
# Row 77
$ws.Range("H77").Value = 2083.8125
$ws.Range("I77").Value = 1742.7
$ws.Range("J77").Value = 2652.3333
$ws.Range("K77").Value = 8713.5
$ws.Range("L77").Value = 13261.6665
$ws.Range("M77").Value = -4345.5
$ws.Range("N77").Value = -21997.6665
# Row 141
$ws.Range("H141").Value = 98000
$ws.Range("J141").Value = 98000
$ws.Range("L141").Value = 98000
$ws.Range("N141").Value = -108360

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -127

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1653.5
$ws.Range("I16").Value = 1575
$ws.Range("J16").Value = 1810.5
$ws.Range("K16").Value = 1575
$ws.Range("L16").Value = 1810.5
$ws.Range("M16").Value = -1288
$ws.Range("N16").Value = -2384.5
# Row 31
$ws.Range("H31").Value = 33239.76
$ws.Range("I31").Value = 835.8570999999999
$ws.Range("J31").Value = 56704.656
$ws.Range("K31").Value = 835.8570999999999
$ws.Range("L31").Value = 56704.656
$ws.Range("M31").Value = -540.8570999999999
$ws.Range("N31").Value = -57294.656
# Row 34
$ws.Range("H34").Value = 33239.76
$ws.Range("I34").Value = 835.8570999999999
$ws.Range("J34").Value = 56704.656
$ws.Range("K34").Value = 835.8570999999999
$ws.Range("L34").Value = 56704.656
$ws.Range("M34").Value = -633.8570999999999
$ws.Range("N34").Value = -57108.656
# Row 99
$ws.Range("H99").Value = 15827.444
$ws.Range("I99").Value = 5200.8
$ws.Range("J99").Value = 29110.75
$ws.Range("K99").Value = 5200.8
$ws.Range("L99").Value = 29110.75
$ws.Range("M99").Value = -3702.8
$ws.Range("N99").Value = -32106.75
# Row 107
$ws.Range("H107").Value = 714.6667
$ws.Range("I107").Value = 985.5714
$ws.Range("J107").Value = 477.625
$ws.Range("K107").Value = 985.5714
$ws.Range("L107").Value = 477.625
$ws.Range("M107").Value = 934.4286
$ws.Range("N107").Value = -4317.625
# Row 113
$ws.Range("H113").Value = 1653.5
$ws.Range("I113").Value = 1575
$ws.Range("J113").Value = 1810.5
$ws.Range("K113").Value = 1575
$ws.Range("L113").Value = 1810.5
$ws.Range("M113").Value = 595
$ws.Range("N113").Value = -6150.5
# Row 126
$ws.Range("H126").Value = 15827.444
$ws.Range("I126").Value = 5200.8
$ws.Range("J126").Value = 29110.75
$ws.Range("K126").Value = 15602.4
$ws.Range("L126").Value = 87332.25
$ws.Range("M126").Value = -13132.4
$ws.Range("N126").Value = -92272.25

$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 380.125
$ws.Range("I11").Value = 346.5
$ws.Range("J11").Value = 413.75
$ws.Range("K11").Value = 1039.5
$ws.Range("L11").Value = 1241.25
$ws.Range("M11").Value = -899.5
$ws.Range("N11").Value = -1521.25
# Row 22
$ws.Range("H22").Value = 7376.5454
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 7376.5454
$ws.Range("K22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("M22").Value = 22129.6362
$ws.Range("N22").Value = -22467.6362
# Row 27
$ws.Range("H27").Value = 7376.5454
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 7376.5454
$ws.Range("K27").Value = 0
$ws.Range("L27").ClearContents()
$ws.Range("M27").Value = 22129.6362
$ws.Range("N27").Value = -22333.6362
# Row 58
$ws.Range("H58").Value = 2582.8333
$ws.Range("I58").Value = 1499.6666
$ws.Range("J58").Value = 3666
$ws.Range("K58").Value = 4498.9998
$ws.Range("L58").Value = 10998
$ws.Range("M58").Value = -4370.9998
$ws.Range("N58").Value = -11254
# Row 124
$ws.Range("H124").Value = 3448.3333
$ws.Range("J124").Value = 4240
$ws.Range("L124").Value = 12720
$ws.Range("N124").Value = -22540
# Row 129
$ws.Range("H129").Value = 13424440
$ws.Range("I129").Value = 62514396
$ws.Range("J129").Value = 333785
$ws.Range("K129").Value = 187543188
$ws.Range("L129").Value = 1001355
$ws.Range("M129").Value = -187538188
$ws.Range("N129").Value = -1011355
# Row 131
$ws.Range("H131").Value = 7122.9375
$ws.Range("I131").Value = 1281.6666
$ws.Range("J131").Value = 7512.3555
$ws.Range("K131").Value = 3844.9998
$ws.Range("L131").Value = 22537.0665
$ws.Range("M131").Value = 1195.0002
$ws.Range("N131").Value = -32617.0665

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 95245.37
$ws.Range("I70").Value = 185536.73
$ws.Range("K70").Value = 185536.73
$ws.Range("M70").Value = -185266.73
# Row 73
$ws.Range("H73").Value = 95245.37
$ws.Range("I73").Value = 185536.73
$ws.Range("K73").Value = 185536.73
$ws.Range("M73").Value = -184600.73

$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 6943.375
$ws.Range("J54").Value = 6925.2856
$ws.Range("L54").Value = 6925.2856
$ws.Range("N54").Value = -7965.2856
# Row 122
$ws.Range("H122").Value = 1754.0714
$ws.Range("I122").Value = 1442.8572
$ws.Range("J122").Value = 2065.2856
$ws.Range("K122").Value = 4328.571599999999
$ws.Range("L122").Value = 6195.8568
$ws.Range("M122").Value = -1878.571599999999
$ws.Range("N122").Value = -11095.8568
